# Apply the "fixed p2 and p5" corrections described in the commit.
#
# All affected cells store plain numeric-looking values as *text*
# (inlineStr) in the original workbook, so we force text formatting
# before writing them, then restore the default "Normal" style so we
# don't leave stray custom number formats behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Part 2 sheet ---
$wsPart2 = $wb.Worksheets.Item("Part 2")
Set-TextValue $wsPart2.Range("D3") "0.017391304347826087"

# --- Part 5 sheet ---
$wsPart5 = $wb.Worksheets.Item("Part 5")
Set-TextValue $wsPart5.Range("B2") "20.0"
Set-TextValue $wsPart5.Range("D2") "1"
Set-TextValue $wsPart5.Range("B3") "8.0"
Set-TextValue $wsPart5.Range("D3") "0"
Set-TextValue $wsPart5.Range("B4") "7.0"
